$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 359
$ws.Range("F5").Value = 17
$ws.Range("F7").Value = 445
$ws.Range("F9").Value = 177
$ws.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202405/pfNAP8zt1715074013459.jpeg"
$ws.Range("F10").Value = 149
$ws.Range("C11").Value = "杭州·D3动漫游戏嘉年华"
$ws.Range("D11").Value = "德胜东路2539号 梦马汽车小镇"
$ws.Range("E11").Value = "2024.05.25 10:00-05.25 17:00"
$ws.Range("F11").Value = 167
$ws.Range("G11").Value = 50
$ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=84912"
$ws.Range("I11").Value = "//i1.hdslb.com/bfs/openplatform/202405/HXY7pTYI1715053764601.jpeg"
$ws.Range("C12").Value = "杭州·Redamancy动漫游戏嘉年华×运动番全明星"
$ws.Range("D12").Value = "富春路80号(甬江路地铁站A口旁) 杭州全民健身中心"
$ws.Range("E12").Value = "2024.05.25 10:00-05.26 17:00"
$ws.Range("F12").Value = 1042
$ws.Range("G12").Value = 68
$ws.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=84947"
$ws.Range("I12").Value = "//i1.hdslb.com/bfs/openplatform/202404/65Usx6BT1713796309433.jpeg"
$ws.Range("C13").Value = "杭州·动漫视界COS盛典"
$ws.Range("D13").Value = "花蒋路1号 岚图汽车杭州元通全功能用户中心"
$ws.Range("E13").Value = "2024.05.25 09:00-05.25 17:00"
$ws.Range("F13").Value = 3
$ws.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=85106"
$ws.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202404/9kTcunDW1714098692172.jpeg"
$ws.Range("C14").Value = "杭州·原神X星铁X绝区零only"
$ws.Range("D14").Value = "望江东路333号 杭州瑞莱克斯大酒店"
$ws.Range("F14").Value = 262
$ws.Range("G14").Value = 60
$ws.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=82754"
$ws.Range("I14").Value = "//i1.hdslb.com/bfs/openplatform/202403/qA0LNJuF1710234461030.jpeg"
$ws.Range("F15").Value = 181
$ws.Range("F16").Value = 1477
$ws.Range("F17").Value = 545
$ws.Range("F18").Value = 222
$ws.Range("F19").Value = 340
$ws.Range("F21").Value = 804
$ws.Range("F22").Value = 1140
$ws.Range("F24").Value = 1917
$ws.Range("F25").Value = 2636
$ws.Range("F26").Value = 1412
$ws.Range("F28").Value = 27
$ws.Range("F29").Value = 370
$ws.Range("F30").Value = 407
$ws.Range("F31").Value = 1163
$ws.Range("F32").Value = 809
$ws.Range("F33").Value = 1294
$ws.Range("F34").Value = 152
$ws.Range("F37").Value = 567
$ws.Range("F38").Value = 657
$ws.Range("F39").Value = 825
$ws.Range("F40").Value = 352
$ws.Range("F41").Value = 237

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F15").Value = 613

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 359
$ws.Range("F11").Value = 445
$ws.Range("F13").Value = 177
$ws.Range("I13").Value = "//i2.hdslb.com/bfs/openplatform/202405/pfNAP8zt1715074013459.jpeg"
$ws.Range("F15").Value = 149
$ws.Range("C16").Value = "杭州·D3动漫游戏嘉年华"
$ws.Range("D16").Value = "德胜东路2539号 梦马汽车小镇"
$ws.Range("E16").Value = "2024.05.25 10:00-05.25 17:00"
$ws.Range("F16").Value = 167
$ws.Range("G16").Value = 50
$ws.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=84912"
$ws.Range("I16").Value = "//i1.hdslb.com/bfs/openplatform/202405/HXY7pTYI1715053764601.jpeg"
$ws.Range("C17").Value = "杭州·动漫视界COS盛典"
$ws.Range("D17").Value = "花蒋路1号 岚图汽车杭州元通全功能用户中心"
$ws.Range("E17").Value = "2024.05.25 09:00-05.25 17:00"
$ws.Range("F17").Value = 3
$ws.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=85106"
$ws.Range("I17").Value = "//i1.hdslb.com/bfs/openplatform/202404/9kTcunDW1714098692172.jpeg"
$ws.Range("C18").Value = "杭州·原神X星铁X绝区零only"
$ws.Range("D18").Value = "望江东路333号 杭州瑞莱克斯大酒店"
$ws.Range("E18").Value = "2024.05.25 10:00-05.25 17:00"
$ws.Range("F18").Value = 262
$ws.Range("G18").Value = 60
$ws.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=82754"
$ws.Range("I18").Value = "//i1.hdslb.com/bfs/openplatform/202403/qA0LNJuF1710234461030.jpeg"
$ws.Range("C19").Value = "杭州·早鸟5折起·《LALALAND爱乐之城》浪漫主题音乐会"
$ws.Range("D19").Value = "武林路77号 浙江省文化馆小剧场（原群艺馆小剧场）"
$ws.Range("E19").Value = "2024.05.25 19:30-05.25 21:00"
$ws.Range("F19").Value = 9
$ws.Range("G19").Value = 100
$ws.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=84519"
$ws.Range("I19").Value = "//i1.hdslb.com/bfs/openplatform/202404/jJLft5tT1712888683239.jpeg"
$ws.Range("F20").Value = 181
$ws.Range("F21").Value = 1477
$ws.Range("F22").Value = 545
$ws.Range("F23").Value = 222
$ws.Range("F24").Value = 340
$ws.Range("F26").Value = 1140
$ws.Range("F27").Value = 2636
$ws.Range("F29").Value = 1412
$ws.Range("F32").Value = 27
$ws.Range("F34").Value = 370
$ws.Range("F35").Value = 407
$ws.Range("F36").Value = 1163
$ws.Range("F39").Value = 809
$ws.Range("F40").Value = 1294
$ws.Range("F42").Value = 567
$ws.Range("F43").Value = 657
$ws.Range("F44").Value = 825
$ws.Range("F45").Value = 352
$ws.Range("F48").Value = 237
